$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (pushes existing rows 7-22 down to 8-23)
$ws.Rows.Item(7).Insert()

# Copy formatting from row 6 (existing data-row style) onto the newly
# inserted, still-blank row 7 so it keeps the same cell style (s="4")
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)

# Populate the new row 7 with the inserted category
$ws.Range("A7").Value = "K_BEWERTUNGTEC"
$ws.Range("B7").Value = "Erfüllungsgrad"
$ws.Range("C7").Value = "Level of compliance"

# Update row 6 (K_BEWERTUNG) German/English labels
$ws.Range("B6").Value = "Bewertung"
$ws.Range("C6").Value = "Rating"

# Widen column B to fit the new/updated labels (target stored width 24.0234375;
# the ColumnWidth API here quantizes to 1/7-character steps, so 163/7 = 23.2857...
# is the input that rounds to the closest achievable stored width of 24)
$ws.Columns.Item(2).ColumnWidth = 23.285714285714285
